$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2611.125
$ws.Range("I28").Value = 593.05554
$ws.Range("K28").Value = 593.05554
$ws.Range("M28").Value = -108.05554

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1024
$ws.Range("I32").Value = 999
$ws.Range("J32").Value = 1032.3334
$ws.Range("K32").Value = 999
$ws.Range("L32").Value = 1032.3334
$ws.Range("M32").Value = -673
$ws.Range("N32").Value = -1684.3334

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 406.5
$ws.Range("I80").Value = 451.66666
$ws.Range("K80").Value = 1354.99998
$ws.Range("M80").Value = -356.9999800000001

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 406.5
$ws.Range("I83").Value = 451.66666
$ws.Range("K83").Value = 4064.99994
$ws.Range("M83").Value = 927.0000600000003

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4223.1523
$ws.Range("I132").Value = 4407.95
$ws.Range("K132").Value = 13223.85
$ws.Range("M132").Value = -10693.85

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3238.0908
$ws.Range("I138").Value = 3508.3333
$ws.Range("J138").Value = 2991.348
$ws.Range("K138").Value = 10524.9999
$ws.Range("L138").Value = 8974.044
$ws.Range("M138").Value = -5384.999899999999
$ws.Range("N138").Value = -19254.044

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1208.5294
$ws.Range("I2").Value = 545
$ws.Range("J2").Value = 2425
$ws.Range("K2").Value = 545
$ws.Range("L2").Value = 2425
$ws.Range("M2").Value = -432
$ws.Range("N2").Value = -2651

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4827.7295
$ws.Range("I32").Value = 4573.2
$ws.Range("K32").Value = 4573.2
$ws.Range("M32").Value = -4286.2

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 30586.867
$ws.Range("I45").Value = 43394.7
$ws.Range("J45").Value = 4971.2
$ws.Range("K45").Value = 43394.7
$ws.Range("L45").Value = 4971.2
$ws.Range("M45").Value = -43017.7
$ws.Range("N45").Value = -5725.2

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2398.5557
$ws.Range("I61").Value = 652.1818
$ws.Range("K61").Value = 652.1818
$ws.Range("M61").Value = -440.1818

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 216232.97
$ws.Range("I74").Value = 695579.4
$ws.Range("J74").Value = 3190.111
$ws.Range("K74").Value = 695579.4
$ws.Range("L74").Value = 3190.111
$ws.Range("M74").Value = -694705.4
$ws.Range("N74").Value = -4938.111

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 216232.97
$ws.Range("I77").Value = 695579.4
$ws.Range("J77").Value = 3190.111
$ws.Range("K77").Value = 3477897
$ws.Range("L77").Value = 15950.555
$ws.Range("M77").Value = -3473529
$ws.Range("N77").Value = -24686.555

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5131.0713
$ws.Range("I102").Value = 4983.5
$ws.Range("K102").Value = 4983.5
$ws.Range("M102").Value = -3361.5

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1208.5294
$ws.Range("I116").Value = 545
$ws.Range("J116").Value = 2425
$ws.Range("K116").Value = 545
$ws.Range("L116").Value = 2425
$ws.Range("M116").Value = 1749
$ws.Range("N116").Value = -7013

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1421.1666
$ws.Range("I132").Value = 666.5
$ws.Range("K132").Value = 1999.5
$ws.Range("M132").Value = 530.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2398.5557
$ws.Range("I136").Value = 652.1818
$ws.Range("K136").Value = 1956.5454
$ws.Range("M136").Value = 593.4546

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1208.5294
$ws.Range("I3").Value = 545
$ws.Range("J3").Value = 2425
$ws.Range("K3").Value = 545
$ws.Range("L3").Value = 2425
$ws.Range("M3").Value = -431
$ws.Range("N3").Value = -2653

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27782770
$ws.Range("I20").Value = 50006424
$ws.Range("K20").Value = 50006424
$ws.Range("M20").Value = -50006177

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1432.8077
$ws.Range("I107").Value = 1180
$ws.Range("K107").Value = 1180
$ws.Range("M107").Value = 740

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1636.12
$ws.Range("I16").Value = 1472
$ws.Range("J16").Value = 2155.8333
$ws.Range("K16").Value = 1472
$ws.Range("L16").Value = 2155.8333
$ws.Range("M16").Value = -1185
$ws.Range("N16").Value = -2729.8333

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5896.2
$ws.Range("I31").Value = 4321
$ws.Range("J31").Value = 6782.25
$ws.Range("K31").Value = 4321
$ws.Range("L31").Value = 6782.25
$ws.Range("M31").Value = -4026
$ws.Range("N31").Value = -7372.25

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5896.2
$ws.Range("I34").Value = 4321
$ws.Range("J34").Value = 6782.25
$ws.Range("K34").Value = 4321
$ws.Range("L34").Value = 6782.25
$ws.Range("M34").Value = -4119
$ws.Range("N34").Value = -7186.25

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2317.0952
$ws.Range("I58").Value = 1164.875
$ws.Range("K58").Value = 1164.875
$ws.Range("M58").Value = -961.875

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1636.12
$ws.Range("I113").Value = 1472
$ws.Range("J113").Value = 2155.8333
$ws.Range("K113").Value = 1472
$ws.Range("L113").Value = 2155.8333
$ws.Range("M113").Value = 698
$ws.Range("N113").Value = -6495.8333

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2426.2188
$ws.Range("J134").Value = 3324.75
$ws.Range("L134").Value = 9974.25
$ws.Range("N134").Value = -15044.25

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2317.0952
$ws.Range("I136").Value = 1164.875
$ws.Range("K136").Value = 3494.625
$ws.Range("M136").Value = -944.625

# CRP row 137
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 105999
$ws.Range("J137").Value = 105999
$ws.Range("L137").Value = 105999
$ws.Range("M137").Value = -116199

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2621.2222
$ws.Range("I3").Value = 1764.8334
$ws.Range("J3").Value = 4334
$ws.Range("K3").Value = 5294.5002
$ws.Range("L3").Value = 13002
$ws.Range("M3").Value = -5182.5002
$ws.Range("N3").Value = -13226

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4311.375
$ws.Range("I81").Value = 2500
$ws.Range("J81").Value = 4915.1665
$ws.Range("K81").Value = 7500
$ws.Range("L81").Value = 14745.4995
$ws.Range("M81").Value = -6377
$ws.Range("N81").Value = -16991.4995

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 4311.375
$ws.Range("I84").Value = 2500
$ws.Range("J84").Value = 4915.1665
$ws.Range("K84").Value = 22500
$ws.Range("L84").Value = 44236.4985
$ws.Range("M84").Value = -16884
$ws.Range("N84").Value = -55468.4985

# CUL row 108
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()

# CUL row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2874.7144
$ws.Range("J109").Value = 2899.6667
$ws.Range("L109").Value = 8699.000100000001
$ws.Range("N109").Value = -10779.0001

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1790.4445
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 131550.56
$ws.Range("I70").Value = 337325.84
$ws.Range("K70").Value = 337325.84
$ws.Range("M70").Value = -337055.84

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 131550.56
$ws.Range("I73").Value = 337325.84
$ws.Range("K73").Value = 337325.84
$ws.Range("M73").Value = -336389.84

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1985.2222
$ws.Range("I132").Value = 1510.1875
$ws.Range("K132").Value = 4530.5625
$ws.Range("M132").Value = -2000.5625

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4139.3
$ws.Range("I136").Value = 4185.933
$ws.Range("K136").Value = 12557.799
$ws.Range("M136").Value = -10007.799

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 434.58334
$ws.Range("I4").Value = 401
$ws.Range("K4").Value = 401
$ws.Range("M4").Value = -288

# WVR row 75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 47704.332
$ws.Range("J75").Value = 47497.5
$ws.Range("L75").Value = 47497.5
$ws.Range("N75").Value = -49369.5

# WVR row 78
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 47704.332
$ws.Range("J78").Value = 47497.5
$ws.Range("L78").Value = 142492.5
$ws.Range("N78").Value = -151852.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3496.8125
$ws.Range("I132").Value = 4293.8887
$ws.Range("J132").Value = 2472
$ws.Range("K132").Value = 12881.6661
$ws.Range("L132").Value = 7416
$ws.Range("M132").Value = -10351.6661
$ws.Range("N132").Value = -12476

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 125001350
$ws.Range("I136").Value = 250001090
$ws.Range("J136").Value = 1618.75
$ws.Range("K136").Value = 750003270
$ws.Range("L136").Value = 4856.25
$ws.Range("M136").Value = -750000720
$ws.Range("N136").Value = -9956.25
